# Refresh of scraped market-board prices / leve profit columns
# (H:N = currentAveragePrice[, NQ, HQ], LevePriceNQ/HQ, LeveProfitNQ/HQ)
# on the per-class crafting-leve-profit sheets. Values below mirror the
# latest Universalis price pull; only numeric cells H..N are touched.
$wb = $excel.ActiveWorkbook

# --- ALC sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# row 4
$ws.Range("H4").Value = 212
$ws.Range("I4").Value = 212
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 212
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -98
$ws.Range("N4").Value = $null

# row 17
$ws.Range("H17").Value = 1063
$ws.Range("I17").Value = 1063
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 3189
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -3021
$ws.Range("N17").Value = $null

# row 92
$ws.Range("H92").Value = 266.7857
$ws.Range("I92").Value = 299.58334
$ws.Range("J92").Value = 70
$ws.Range("K92").Value = 299.58334
$ws.Range("L92").Value = 70
$ws.Range("M92").Value = 948.41666
$ws.Range("N92").Value = -2566

# row 97
$ws.Range("H97").Value = 500001000
$ws.Range("J97").Value = 500001000
$ws.Range("L97").Value = 1500003000
$ws.Range("N97").Value = -1500003992

# row 98
$ws.Range("H98").Value = 2131.5334
$ws.Range("I98").Value = 696.1818
$ws.Range("J98").Value = 6078.75
$ws.Range("K98").Value = 696.1818
$ws.Range("L98").Value = 6078.75
$ws.Range("M98").Value = 801.8182
$ws.Range("N98").Value = -9074.75

# row 99
$ws.Range("H99").Value = 1579
$ws.Range("I99").Value = 1678.5
$ws.Range("J99").Value = 1181
$ws.Range("K99").Value = 5035.5
$ws.Range("L99").Value = 3543
$ws.Range("M99").Value = -3537.5
$ws.Range("N99").Value = -6539

# row 122
$ws.Range("H122").Value = 2131.5334
$ws.Range("I122").Value = 696.1818
$ws.Range("J122").Value = 6078.75
$ws.Range("K122").Value = 2088.5454
$ws.Range("L122").Value = 18236.25
$ws.Range("M122").Value = 361.4546
$ws.Range("N122").Value = -23136.25

# row 125
$ws.Range("H125").Value = 20836114
$ws.Range("I125").Value = 31251540
$ws.Range("K125").Value = 281263860
$ws.Range("M125").Value = -281261400

# row 127
$ws.Range("H127").Value = 694.5
$ws.Range("I127").Value = 694.5
$ws.Range("K127").Value = 2083.5
$ws.Range("M127").Value = 2876.5

# row 137
$ws.Range("H137").Value = 5277.5
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = $null


# --- ARM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4918.0464
$ws.Range("I32").Value = 3694.561
$ws.Range("K32").Value = 3694.561
$ws.Range("M32").Value = -3407.561

# row 61
$ws.Range("H61").Value = 7999
$ws.Range("J61").Value = 7999
$ws.Range("L61").Value = 7999
$ws.Range("N61").Value = -8423

# row 74
$ws.Range("H74").Value = 2849
$ws.Range("I74").Value = 2849
$ws.Range("K74").Value = 2849
$ws.Range("M74").Value = -1975

# row 77
$ws.Range("H77").Value = 2849
$ws.Range("I77").Value = 2849
$ws.Range("K77").Value = 14245
$ws.Range("M77").Value = -9877

# row 97
$ws.Range("H97").Value = 459.33334
$ws.Range("I97").Value = 435.5
$ws.Range("J97").Value = 650
$ws.Range("K97").Value = 435.5
$ws.Range("L97").Value = 650
$ws.Range("M97").Value = 60.5
$ws.Range("N97").Value = -1642

# row 102
$ws.Range("H102").Value = 1938.0769
$ws.Range("I102").Value = 1199.5454
$ws.Range("K102").Value = 1199.5454
$ws.Range("M102").Value = 422.4546

# row 132
$ws.Range("H132").Value = 1746
$ws.Range("I132").Value = 1746
$ws.Range("K132").Value = 5238
$ws.Range("M132").Value = -2708

# row 136
$ws.Range("H136").Value = 7999
$ws.Range("J136").Value = 7999
$ws.Range("L136").Value = 23997
$ws.Range("N136").Value = -29097


# --- BSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# row 61
$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25626

# row 99
$ws.Range("H99").Value = 1376.125
$ws.Range("I99").Value = 1251.5
$ws.Range("K99").Value = 1251.5
$ws.Range("M99").Value = 246.5

# row 107
$ws.Range("H107").Value = 1504
$ws.Range("I107").Value = 1448.8889
$ws.Range("K107").Value = 1448.8889
$ws.Range("M107").Value = 471.1111000000001

# row 134
$ws.Range("H134").Value = 4997
$ws.Range("I134").Value = 4997
$ws.Range("K134").Value = 14991
$ws.Range("M134").Value = -12456


# --- CRP sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# row 8
$ws.Range("H8").Value = 800
$ws.Range("I8").Value = 800
$ws.Range("J8").Value = 800
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = -660
$ws.Range("N8").Value = -1080

# row 16
$ws.Range("H16").Value = 2453.8572
$ws.Range("I16").Value = 700
$ws.Range("J16").Value = 3155.4
$ws.Range("K16").Value = 700
$ws.Range("L16").Value = 3155.4
$ws.Range("M16").Value = -413
$ws.Range("N16").Value = -3729.4

# row 22
$ws.Range("H22").Value = 1565.1111
$ws.Range("I22").Value = 1222.5
$ws.Range("J22").Value = 1839.2
$ws.Range("K22").Value = 1222.5
$ws.Range("L22").Value = 1839.2
$ws.Range("M22").Value = -872.5
$ws.Range("N22").Value = -2539.2

# row 58
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null

# row 62
$ws.Range("H62").Value = 3085
$ws.Range("I62").Value = 3481.25
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 3481.25
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -2857.25
$ws.Range("N62").Value = -2748

# row 65
$ws.Range("H65").Value = 3085
$ws.Range("I65").Value = 3481.25
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 17406.25
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -14286.25
$ws.Range("N65").Value = -13740

# row 93
$ws.Range("H93").Value = 24999.5
$ws.Range("I93").Value = 24999.5
$ws.Range("K93").Value = 24999.5
$ws.Range("M93").Value = -23127.5

# row 113
$ws.Range("H113").Value = 2453.8572
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 3155.4
$ws.Range("K113").Value = 700
$ws.Range("L113").Value = 3155.4
$ws.Range("M113").Value = 1470
$ws.Range("N113").Value = -7495.4

# row 132
$ws.Range("H132").Value = 10612.059
$ws.Range("I132").Value = 4343.7144
$ws.Range("J132").Value = 14999.9
$ws.Range("K132").Value = 13031.1432
$ws.Range("L132").Value = 44999.7
$ws.Range("M132").Value = -10501.1432
$ws.Range("N132").Value = -50059.7

# row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

# row 100
$ws.Range("H100").Value = 14000
$ws.Range("J100").Value = 14000
$ws.Range("L100").Value = 42000
$ws.Range("N100").Value = -43622


# --- CUL sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# row 117
$ws.Range("H117").Value = 150
$ws.Range("I117").Value = 150
$ws.Range("K117").Value = 450
$ws.Range("M117").Value = 2992

# row 140
$ws.Range("H140").Value = 1969.6666
$ws.Range("I140").Value = 704.5
$ws.Range("K140").Value = 2113.5
$ws.Range("M140").Value = 3066.5


# --- GSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# row 11
$ws.Range("H11").Value = 801960.6
$ws.Range("I11").Value = 1336067.6
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 1336067.6
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = -1335928.6
$ws.Range("N11").Value = -1078

# row 80
$ws.Range("H80").Value = 2811.5
$ws.Range("J80").Value = 2728.6
$ws.Range("L80").Value = 2728.6
$ws.Range("N80").Value = -4724.6

# row 83
$ws.Range("H83").Value = 2811.5
$ws.Range("J83").Value = 2728.6
$ws.Range("L83").Value = 13643
$ws.Range("N83").Value = -23627

# row 122
$ws.Range("H122").Value = 3831.5
$ws.Range("I122").Value = 3799.7334
$ws.Range("K122").Value = 11399.2002
$ws.Range("M122").Value = -8949.200199999999

# row 132
$ws.Range("H132").Value = 3561.75
$ws.Range("I132").Value = 3067.3635
$ws.Range("K132").Value = 9202.0905
$ws.Range("M132").Value = -6672.0905


# --- LTW sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2182.6
$ws.Range("I7").Value = 2134
$ws.Range("J7").Value = 2296
$ws.Range("K7").Value = 2134
$ws.Range("L7").Value = 2296
$ws.Range("M7").Value = -2022
$ws.Range("N7").Value = -2520

# row 29
$ws.Range("H29").Value = 14333.333
$ws.Range("I29").Value = 15500
$ws.Range("K29").Value = 15500
$ws.Range("M29").Value = -15205

# row 46
$ws.Range("H46").Value = 1506.6923
$ws.Range("I46").Value = 1580.75
$ws.Range("K46").Value = 1580.75
$ws.Range("M46").Value = -1392.75

# row 82
$ws.Range("H82").Value = 3194.4
$ws.Range("I82").Value = 2314.6667
$ws.Range("J82").Value = 3571.4285
$ws.Range("K82").Value = 2314.6667
$ws.Range("L82").Value = 3571.4285
$ws.Range("M82").Value = -1953.6667
$ws.Range("N82").Value = -4293.4285

# row 85
$ws.Range("H85").Value = 3194.4
$ws.Range("I85").Value = 2314.6667
$ws.Range("J85").Value = 3571.4285
$ws.Range("K85").Value = 2314.6667
$ws.Range("L85").Value = 3571.4285
$ws.Range("M85").Value = -1066.6667
$ws.Range("N85").Value = -6067.4285

# row 122
$ws.Range("H122").Value = 8990
$ws.Range("I122").Value = 8990
$ws.Range("K122").Value = 26970
$ws.Range("M122").Value = -24520

# row 126
$ws.Range("H126").Value = 2182.6
$ws.Range("I126").Value = 2134
$ws.Range("J126").Value = 2296
$ws.Range("K126").Value = 6402
$ws.Range("L126").Value = 6888
$ws.Range("M126").Value = -3932
$ws.Range("N126").Value = -11828


# --- WVR sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null

# row 96
$ws.Range("H96").Value = 3670
$ws.Range("I96").Value = 3466.6667
$ws.Range("J96").Value = 5500
$ws.Range("K96").Value = 3466.6667
$ws.Range("L96").Value = 5500
$ws.Range("M96").Value = -2093.6667
$ws.Range("N96").Value = -8246

# row 132
$ws.Range("H132").Value = 2259.32
$ws.Range("I132").Value = 2126.125
$ws.Range("K132").Value = 6378.375
$ws.Range("M132").Value = -3848.375
